# Update the "cvd" (column E) values for the Manufacturing Voluntary
# Turnover rows across every site worksheet: the CVD figure moved from
# 0.0639 to 0.0776. Walk every worksheet instead of hard-coding a sheet
# list / cell list so this keeps working if more site tabs are added
# later (matches the intent of "need to make that dynamic").

$wb = $excel.ActiveWorkbook

$oldCvd = 0.0639
$newCvd = 0.0776
$tolerance = 0.0000001
$cvdColumn = 5   # column E

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        $cell = $ws.Cells.Item($r, $cvdColumn)
        $val = $cell.Value2
        if ($val -ne $null -and [math]::Abs($val - $oldCvd) -lt $tolerance) {
            $cell.Value2 = $newCvd
        }
    }
}
